# Implement first version of lot sizing rules:
# - Increase NrBuckets (Generic!B4) from 4 to 5
# - Extend ForecastedAverageDemand and ForcastedStandardDeviation sheets
#   with a new row (row 6) that repeats the pattern of the last existing
#   row (row 5), with the bucket index in column A incremented.

$wb = $excel.ActiveWorkbook

# 1. Update the Generic sheet's NrBuckets value (B4: 4 -> 5)
$wsGeneric = $wb.Worksheets.Item("Generic")
$wsGeneric.Range("B4").Value = 5

# 2. Append the new bucket row to both demand-related sheets
$sheetNames = @("ForecastedAverageDemand", "ForcastedStandardDeviation")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Copy row 5 (A5:W5) down into the new row 6 (A6:W6), including formatting
    $ws.Range("A5:W5").Copy($ws.Range("A6:W6"))

    # Column A holds the zero-based bucket index; increment it for the new row
    $lastIndex = [double]$ws.Range("A5").Value2
    $ws.Range("A6").Value = $lastIndex + 1
}
